# Commit: "IWP TestData, Demo Verification Script: Committing after Windows 11 update"
# The Katalon demo/test-automation suite was re-run after the Windows 11 update,
# which refreshed the Result/Date columns on the log sheets touched by that run.

$wb = $excel.ActiveWorkbook

# --- PayNowCC: rows 2-7, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Wed Jun 25 00:45:49 IST 2025"
$ws.Range("B3").Value = "Wed Jun 25 00:46:35 IST 2025"
$ws.Range("B4").Value = "Wed Jun 25 00:47:23 IST 2025"
$ws.Range("B5").Value = "Wed Jun 25 00:48:14 IST 2025"
$ws.Range("B6").Value = "Wed Jun 25 00:49:03 IST 2025"
$ws.Range("B7").Value = "Wed Jun 25 00:49:51 IST 2025"

# --- NoModifyAmountCC: row 2, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("B2").Value = "Wed Jun 25 00:35:23 IST 2025"

# --- NoModifyBillingAddressCC: row 2, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("NoModifyBillingAddressCC")
$ws.Range("B2").Value = "Wed Jun 25 00:39:06 IST 2025"

# --- CMCAutopayCC: row 2, Result flips Fail -> Pass, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("CMCAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Jun 25 00:29:11 IST 2025"

# --- PayNowCreditCardDCF: row 2, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("PayNowCreditCardDCF")
$ws.Range("B2").Value = "Wed Jun 25 00:42:31 IST 2025"

# --- PayNowCreditCardSCF: row 2, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("PayNowCreditCardSCF")
$ws.Range("B2").Value = "Wed Jun 25 00:52:51 IST 2025"

# --- DCFCCVerbiage: row 2, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("DCFCCVerbiage")
$ws.Range("B2").Value = "Wed Jun 25 00:56:07 IST 2025"

# --- SCFCCVerbiage: row 2, Date column (B) refreshed ---
$ws = $wb.Worksheets.Item("SCFCCVerbiage")
$ws.Range("B2").Value = "Wed Jun 25 00:58:07 IST 2025"
